$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look numeric,
# so Excel stores them as literal text (matching the source data) instead
# of auto-converting them to numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "68.671.90"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.432.61"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "559.11"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "160.77"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "0.165"
$ws.Range("E9").Value = "  +9.84%  "
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").Value = "0.331"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D13").Value = "68.547.23"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000174"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.873.93"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "23.15"
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("D17").Value = "2.428.54"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "10.49"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").Value = "335.88"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "6.90"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "3.83"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "66.82"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "3.68"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "2.553.67"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "8.22"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "0.0₃0821"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").Value = "7.15"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "428.47"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "1.15"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "160.72"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").Value = "19.07"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D38").Value = "17.92"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("D40").Value = "0.298"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "4.34"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "131.49"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "0.0714"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "0.482"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("E51").Value = "  +0.27%  "
